$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 2025 round 8 results: update Career 180s totals and darts used

# Máté Vass (row 3): Career 180s 1 -> 2
$ws.Range("H3").Value = 2

# Laci Ferenczi (row 5): Darts Used changed, Career 180s 7 -> 8
$ws.Range("G5").Value = "20g One80 Dragon"
$ws.Range("H5").Value = 8

# Ármin Szücs (row 10): Career 180s 0 -> 1
$ws.Range("H10").Value = 1

# Dani Boldizsár (row 16): Career 180s 1 -> 2
$ws.Range("H16").Value = 2

# Update the active selection to reflect the new working cell
$ws.Range("H18").Select()

$wb.Save()
